# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newer snapshot of the data (gh-pages output regenerated
# at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9965
$ws1.Range("F9").Value = 741
$ws1.Range("F11").Value = 1224
$ws1.Range("F13").Value = 3090
$ws1.Range("F14").Value = 2329
$ws1.Range("F16").Value = 2033
$ws1.Range("F17").Value = 244
$ws1.Range("F18").Value = 1939
$ws1.Range("F20").Value = 1578
$ws1.Range("F21").Value = 536
$ws1.Range("F22").Value = 45
$ws1.Range("F25").Value = 41
$ws1.Range("F26").Value = 359
$ws1.Range("F28").Value = 348
$ws1.Range("F29").Value = 559
$ws1.Range("F31").Value = 210
$ws1.Range("F33").Value = 285
$ws1.Range("F34").Value = 1617
$ws1.Range("F35").Value = 89
$ws1.Range("F36").Value = 397
$ws1.Range("F38").Value = 425
$ws1.Range("F39").Value = 899
$ws1.Range("F41").Value = 337

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9965
$ws4.Range("F11").Value = 741
$ws4.Range("F13").Value = 1224
$ws4.Range("F15").Value = 3090
$ws4.Range("F16").Value = 2329
$ws4.Range("F18").Value = 2033
$ws4.Range("F19").Value = 244
$ws4.Range("F20").Value = 1939
$ws4.Range("F22").Value = 1578
$ws4.Range("F23").Value = 536
$ws4.Range("F24").Value = 45
$ws4.Range("F27").Value = 41
$ws4.Range("F28").Value = 359
$ws4.Range("F30").Value = 348
$ws4.Range("F31").Value = 559
$ws4.Range("F36").Value = 210
$ws4.Range("F39").Value = 285
$ws4.Range("F40").Value = 1617
$ws4.Range("F41").Value = 89
$ws4.Range("F43").Value = 397
$ws4.Range("F45").Value = 425
$ws4.Range("F46").Value = 899
$ws4.Range("F48").Value = 337

$wb.Save()
